$wb = $excel.ActiveWorkbook

# ===== Sheet ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1322.2941
$ws.Range("I18").Value = 1217.4375
$ws.Range("K18").Value = 1217.4375
$ws.Range("M18").Value = -933.4375

# Row 33
$ws.Range("H33").Value = 350.07144
$ws.Range("I33").Value = 367.80768
$ws.Range("J33").Value = 119.5
$ws.Range("K33").Value = 367.80768
$ws.Range("L33").Value = 119.5
$ws.Range("M33").Value = -138.80768
$ws.Range("N33").Value = -577.5

# Row 74
$ws.Range("H74").Value = 4213.8887
$ws.Range("I74").Value = 4650
$ws.Range("J74").Value = 3936.3635
$ws.Range("K74").Value = 4650
$ws.Range("L74").Value = 3936.3635
$ws.Range("M74").Value = -3714
$ws.Range("N74").Value = -5808.363499999999

# Row 77
$ws.Range("H77").Value = 4213.8887
$ws.Range("I77").Value = 4650
$ws.Range("J77").Value = 3936.3635
$ws.Range("K77").Value = 23250
$ws.Range("L77").Value = 19681.8175
$ws.Range("M77").Value = -18570
$ws.Range("N77").Value = -29041.8175

# Row 100
$ws.Range("H100").Value = 1250.5555
$ws.Range("I100").Value = 1156.875
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1156.875
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -615.875
$ws.Range("N100").Value = -3082

# Row 132
$ws.Range("H132").Value = 3665079
$ws.Range("I132").Value = 4083562.5
$ws.Range("J132").Value = 3349.75
$ws.Range("K132").Value = 12250687.5
$ws.Range("L132").Value = 10049.25
$ws.Range("M132").Value = -12248157.5
$ws.Range("N132").Value = -15109.25

# Row 137
$ws.Range("H137").Value = 1325.8667
$ws.Range("I137").Value = 893.8421
$ws.Range("K137").Value = 2681.5263
$ws.Range("M137").Value = -131.5263

# Row 138
$ws.Range("H138").Value = 7247539.5
$ws.Range("I138").Value = 9524549
$ws.Range("J138").Value = 2510.6365
$ws.Range("K138").Value = 28573647
$ws.Range("L138").Value = 7531.9095
$ws.Range("M138").Value = -28568507
$ws.Range("N138").Value = -17811.9095

# ===== Sheet ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 102
$ws.Range("H102").Value = 1523.1111
$ws.Range("I102").Value = 1386.8572
$ws.Range("K102").Value = 1386.8572
$ws.Range("M102").Value = 235.1428000000001

# Row 133
$ws.Range("H133").Value = 60000
$ws.Range("J133").Value = 60000
$ws.Range("L133").Value = 60000
$ws.Range("N133").Value = -65060

# ===== Sheet BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2838.375
$ws.Range("I86").Value = 1874.8667
$ws.Range("K86").Value = 1874.8667
$ws.Range("M86").Value = -751.8667

# Row 89
$ws.Range("H89").Value = 2838.375
$ws.Range("I89").Value = 1874.8667
$ws.Range("K89").Value = 9374.333500000001
$ws.Range("M89").Value = -3758.333500000001

# Row 94
$ws.Range("H94").Value = 719.7692
$ws.Range("I94").Value = 704.75
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 704.75
$ws.Range("L94").Value = 900
$ws.Range("M94").Value = -253.75
$ws.Range("N94").Value = -1802

# Row 99
$ws.Range("H99").Value = 722.5833
$ws.Range("I99").Value = 614.5454999999999
$ws.Range("J99").Value = 1911
$ws.Range("K99").Value = 614.5454999999999
$ws.Range("L99").Value = 1911
$ws.Range("M99").Value = 883.4545000000001
$ws.Range("N99").Value = -4907

# Row 105
$ws.Range("H105").Value = 2493.2424
$ws.Range("I105").Value = 2284.3333
$ws.Range("J105").Value = 3433.3333
$ws.Range("K105").Value = 2284.3333
$ws.Range("L105").Value = 3433.3333
$ws.Range("M105").Value = -537.3332999999998
$ws.Range("N105").Value = -6927.3333

# ===== Sheet CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 459.16666
$ws.Range("I22").Value = 383.33334
$ws.Range("J22").Value = 610.8333
$ws.Range("K22").Value = 383.33334
$ws.Range("L22").Value = 610.8333
$ws.Range("M22").Value = -33.33334000000002
$ws.Range("N22").Value = -1310.8333

# Row 31
$ws.Range("H31").Value = 4169086.2
$ws.Range("I31").Value = 1681.3
$ws.Range("K31").Value = 1681.3
$ws.Range("M31").Value = -1386.3

# Row 34
$ws.Range("H34").Value = 4169086.2
$ws.Range("I34").Value = 1681.3
$ws.Range("K34").Value = 1681.3
$ws.Range("M34").Value = -1479.3

# Row 60
$ws.Range("H60").Value = 14719.25
$ws.Range("J60").Value = 14719.25
$ws.Range("L60").Value = 14719.25
$ws.Range("N60").Value = -15741.25

# Row 99
$ws.Range("H99").Value = 49707.855
$ws.Range("I99").Value = 101790.3
$ws.Range("J99").Value = 2360.182
$ws.Range("K99").Value = 101790.3
$ws.Range("L99").Value = 2360.182
$ws.Range("M99").Value = -100292.3
$ws.Range("N99").Value = -5356.182

# Row 105
$ws.Range("H105").Value = 1331.6154
$ws.Range("I105").Value = 1057.1428
$ws.Range("J105").Value = 1651.8334
$ws.Range("K105").Value = 1057.1428
$ws.Range("L105").Value = 1651.8334
$ws.Range("M105").Value = 689.8571999999999
$ws.Range("N105").Value = -5145.8334

# Row 126
$ws.Range("H126").Value = 49707.855
$ws.Range("I126").Value = 101790.3
$ws.Range("J126").Value = 2360.182
$ws.Range("K126").Value = 305370.9
$ws.Range("L126").Value = 7080.545999999999
$ws.Range("M126").Value = -302900.9
$ws.Range("N126").Value = -12020.546

# Row 132
$ws.Range("H132").Value = 2630.879
$ws.Range("I132").Value = 1819.0454
$ws.Range("J132").Value = 4254.5454
$ws.Range("K132").Value = 5457.1362
$ws.Range("L132").Value = 12763.6362
$ws.Range("M132").Value = -2927.1362
$ws.Range("N132").Value = -17823.6362

# ===== Sheet CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 98
$ws.Range("H98").Value = 100.5
$ws.Range("I98").Value = 100
$ws.Range("J98").Value = 100.666664
$ws.Range("K98").Value = 300
$ws.Range("L98").Value = 301.999992
$ws.Range("M98").Value = 1198
$ws.Range("N98").Value = -3297.999992

# Row 131
$ws.Range("H131").Value = 332907.16
$ws.Range("I131").Value = 4435.926
$ws.Range("J131").Value = 529989.9
$ws.Range("K131").Value = 13307.778
$ws.Range("L131").Value = 1589969.7
$ws.Range("M131").Value = -8267.778000000002
$ws.Range("N131").Value = -1600049.7

# ===== Sheet GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 664.44446
$ws.Range("I97").Value = 575
$ws.Range("J97").Value = 736
$ws.Range("K97").Value = 575
$ws.Range("L97").Value = 736
$ws.Range("M97").Value = -79
$ws.Range("N97").Value = -1728

# Row 113
$ws.Range("H113").Value = 35715296
$ws.Range("I113").Value = 250000000
$ws.Range("J113").Value = 1178.1666
$ws.Range("K113").Value = 250000000
$ws.Range("L113").Value = 1178.1666
$ws.Range("M113").Value = -249997830
$ws.Range("N113").Value = -5518.1666

# Row 140
$ws.Range("H140").Value = 44500
$ws.Range("J140").Value = 44500
$ws.Range("L140").Value = 44500
$ws.Range("N140").Value = -54860

# ===== Sheet LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3145
$ws.Range("I22").Value = 3993.3333
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 3993.3333
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -3698.3333
$ws.Range("N22").Value = -1190

# Row 27
$ws.Range("H27").Value = 3145
$ws.Range("I27").Value = 3993.3333
$ws.Range("J27").Value = 600
$ws.Range("K27").Value = 3993.3333
$ws.Range("L27").Value = 600
$ws.Range("M27").Value = -3886.3333
$ws.Range("N27").Value = -814

# Row 62
$ws.Range("H62").Value = 46000
$ws.Range("I62").Value = 46000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 46000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -45376
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 46000
$ws.Range("I65").Value = 46000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 138000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -134880
$ws.Range("N65").ClearContents()

# ===== Sheet WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2660.6
$ws.Range("I62").Value = 2766.6667
$ws.Range("J62").Value = 2501.5
$ws.Range("K62").Value = 2766.6667
$ws.Range("L62").Value = 2501.5
$ws.Range("M62").Value = -2142.6667
$ws.Range("N62").Value = -3749.5

# Row 65
$ws.Range("H65").Value = 2660.6
$ws.Range("I65").Value = 2766.6667
$ws.Range("J65").Value = 2501.5
$ws.Range("K65").Value = 13833.3335
$ws.Range("L65").Value = 12507.5
$ws.Range("M65").Value = -10713.3335
$ws.Range("N65").Value = -18747.5

# Row 96
$ws.Range("H96").Value = 83334000

# Row 113
$ws.Range("H113").Value = 1496.6666
$ws.Range("I113").Value = 1536
$ws.Range("J113").Value = 1300
$ws.Range("K113").Value = 4608
$ws.Range("L113").Value = 3900
$ws.Range("M113").Value = -8240
$ws.Range("N113").Value = -8240
